$d = $word.ActiveDocument

function Set-ParagraphXml($paraRange, [string]$innerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $paraRange.InsertXML($pkg)
}

$tbl = $d.Tables.Item(2)

# Row 2 ("Truong Gia Phu") - status cell: "To Do" -> "Done"
$p1 = $tbl.Cell(2, 6).Range.Paragraphs(1).Range
$xml1 = '<w:p w14:paraId="4BC410BF" w14:textId="77777777" w:rsidR="00964C20" w:rsidRPr="00964C20" w:rsidRDefault="00964C20" w:rsidP="00964C20">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Done</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p1 $xml1

# Row 3 ("Huynh Phat") - status cell: "To Do" -> "Done"
$p2 = $tbl.Cell(3, 6).Range.Paragraphs(1).Range
$xml2 = '<w:p w14:paraId="0417CA48" w14:textId="77777777" w:rsidR="00964C20" w:rsidRPr="00964C20" w:rsidRDefault="00964C20" w:rsidP="00964C20">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Done</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p2 $xml2

# Row 4 ("Nguyen Thanh Dat") - status cell: "To" + " do" (two runs) -> single run "Done"
$p3 = $tbl.Cell(4, 6).Range.Paragraphs(1).Range
$xml3 = '<w:p w14:paraId="3D64F08D" w14:textId="51EBF27F" w:rsidR="00C254B6" w:rsidRPr="00F41917" w:rsidRDefault="00F41917" w:rsidP="00964C20">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Done</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p3 $xml3

# Row 5 ("Nguyen Van Sang Hen") - status cell is empty -> add run "Done" (no lang)
$p4 = $tbl.Cell(5, 6).Range.Paragraphs(1).Range
$xml4 = '<w:p w14:paraId="7BA1B4DC" w14:textId="77777777" w:rsidR="00F41917" w:rsidRDefault="00F41917" w:rsidP="00964C20">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Done</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p4 $xml4

Write-Output "Applied Sprint Planning status updates"
